$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.142470422916503
$ws.Range("D2").Value = 0.1662913088126157
$ws.Range("E2").Value = 0.132469205662674
$ws.Range("F2").Value = 1.387839352556867
$ws.Range("G2").Value = 0.002423862616003382
$ws.Range("J2").Value = 0.139352625704575
$ws.Range("M2").Value = 0.8149744815221283
$ws.Range("N2").Value = 1.527816983927778
$ws.Range("O2").Value = 3.488605281458717
$ws.Range("B3").Value = 0.1329216869159495
$ws.Range("D3").Value = 0.1675120863414179
$ws.Range("E3").Value = 0.1341961599465806
$ws.Range("F3").Value = 1.366939048012213
$ws.Range("G3").Value = 0.002427893894604665
$ws.Range("J3").Value = 0.1417655501149913
$ws.Range("M3").Value = 0.7363472168535452
$ws.Range("N3").Value = 1.476584459131715
$ws.Range("O3").Value = 3.411136809731317
$ws.Range("B4").Value = 0.1271248160866492
$ws.Range("D4").Value = 0.1683116973627001
$ws.Range("E4").Value = 0.1353204602427773
$ws.Range("F4").Value = 1.355023274115169
$ws.Range("G4").Value = 0.002430502538738944
$ws.Range("J4").Value = 0.1433348768683231
$ws.Range("M4").Value = 0.6880271958531523
$ws.Range("N4").Value = 1.445668668305274
$ws.Range("O4").Value = 3.36594553770118
$ws.Range("B5").Value = 0.1247793806297466
$ws.Range("D5").Value = 0.1686501441887067
$ws.Range("E5").Value = 0.1357946937302366
$ws.Range("F5").Value = 1.350397555799987
$ws.Range("G5").Value = 0.002431599240941079
$ws.Range("J5").Value = 0.1439964284945674
$ws.Range("M5").Value = 0.6683270962711418
$ws.Range("N5").Value = 1.433208240783443
$ws.Range("O5").Value = 3.348125074287651
$ws.Range("B6").Value = 0.1243909461329906
$ws.Range("D6").Value = 0.1687071043687673
$ws.Range("E6").Value = 0.1358744103417338
$ws.Range("F6").Value = 1.349643337511935
$ws.Range("G6").Value = 0.002431783383639754
$ws.Range("J6").Value = 0.1441076086669684
$ws.Range("M6").Value = 0.6650553895840261
$ws.Range("N6").Value = 1.431147588498618
$ws.Range("O6").Value = 3.345201905520781
$ws.Range("B7").Value = 0.1270931162646178
$ws.Range("D7").Value = 0.1683162107433027
$ws.Range("E7").Value = 0.1353267908560243
$ws.Range("F7").Value = 1.35495995931187
$ws.Range("G7").Value = 0.002430517192888657
$ws.Range("J7").Value = 0.1433437095998036
$ws.Range("M7").Value = 0.687761549206968
$ws.Range("N7").Value = 1.445500061270366
$ws.Range("O7").Value = 3.36570279654282
$ws.Range("B8").Value = 0.1391644283687299
$ws.Range("D8").Value = 0.1667018535863267
$ws.Range("E8").Value = 0.1330513839458258
$ws.Range("F8").Value = 1.380442238126989
$ws.Range("G8").Value = 0.002425224974357891
$ws.Range("J8").Value = 0.14016636214303
$ws.Range("M8").Value = 0.7878734504521105
$ws.Range("N8").Value = 1.510040814827732
$ws.Range("O8").Value = 3.461400174818834
$ws.Range("B9").Value = 0.1633523989194572
$ws.Range("D9").Value = 0.163932588936321
$ws.Range("E9").Value = 0.1290969427694142
$ws.Range("F9").Value = 1.437719429637781
$ws.Range("G9").Value = 0.002415900661225448
$ws.Range("J9").Value = 0.1346335827846277
$ws.Range("M9").Value = 0.9838029102138535
$ws.Range("N9").Value = 1.64082279609093
$ws.Range("O9").Value = 3.668000021716864
$ws.Range("B10").Value = 0.1814287983962402
$ws.Range("D10").Value = 0.1621388320639685
$ws.Range("E10").Value = 0.1265013585076824
$ws.Range("F10").Value = 1.484303234092607
$ws.Range("G10").Value = 0.002409685550447785
$ws.Range("J10").Value = 0.1309960710820057
$ws.Range("M10").Value = 1.127459237178797
$ws.Range("N10").Value = 1.739388965418186
$ws.Range("O10").Value = 3.831487616063839
$ws.Range("B11").Value = 0.1897167264876316
$ws.Range("D11").Value = 0.1613749107968356
$ws.Range("E11").Value = 0.1253878781364097
$ws.Range("F11").Value = 1.50648386931995
$ws.Range("G11").Value = 0.002406994648712165
$ws.Range("J11").Value = 0.129434491819528
$ws.Range("M11").Value = 1.192737807568136
$ws.Range("N11").Value = 1.784749715473311
$ws.Range("O11").Value = 3.908437286722858
$ws.Range("B12").Value = 0.1928642832981353
$ws.Range("D12").Value = 0.1610931076101103
$ws.Range("E12").Value = 0.1249759128202976
$ws.Range("F12").Value = 1.515026127442241
$ws.Range("G12").Value = 0.002405995173706762
$ws.Range("J12").Value = 0.1288565947248816
$ws.Range("M12").Value = 1.217445624785967
$ws.Range("N12").Value = 1.802000057820237
$ws.Range("O12").Value = 3.937949278330905
$ws.Range("B13").Value = 0.1921859996620441
$ws.Range("D13").Value = 0.1611534665761969
$ws.Range("E13").Value = 0.1250642059459182
$ws.Range("F13").Value = 1.513180031850837
$ws.Range("G13").Value = 0.002406209562352974
$ws.Range("J13").Value = 0.1289804567747908
$ws.Range("M13").Value = 1.212124900971318
$ws.Range("N13").Value = 1.798281662320221
$ws.Range("O13").Value = 3.93157672859337
$ws.Range("B14").Value = 0.1899754968646619
$ws.Range("D14").Value = 0.1613515768816018
$ws.Range("E14").Value = 0.1253537913141036
$ws.Range("F14").Value = 1.507183777679685
$ws.Range("G14").Value = 0.002406912030640273
$ws.Range("J14").Value = 0.1293866782787862
$ws.Range("M14").Value = 1.194770778770973
$ws.Range("N14").Value = 1.786167455223449
$ws.Range("O14").Value = 3.910857772284714
$ws.Range("B15").Value = 0.1886226780633109
$ws.Range("D15").Value = 0.1614738986218782
$ws.Range("E15").Value = 0.1255324324460236
$ws.Range("F15").Value = 1.503529532466615
$ws.Range("G15").Value = 0.002407344850813398
$ws.Range("J15").Value = 0.1296372521908893
$ws.Range("M15").Value = 1.184139302147443
$ws.Range("N15").Value = 1.778756630995275
$ws.Range("O15").Value = 3.898215427041293
$ws.Range("B16").Value = 0.1808884486025306
$ws.Range("D16").Value = 0.1621898029878786
$ws.Range("E16").Value = 0.1265754817046822
$ws.Range("F16").Value = 1.482873643052159
$ws.Range("G16").Value = 0.002409864142871009
$ws.Range("J16").Value = 0.131100002096197
$ws.Range("M16").Value = 1.12319156434377
$ws.Range("N16").Value = 1.736434882872828
$ws.Range("O16").Value = 3.826510832564679
$ws.Range("B17").Value = 0.176160202532202
$ws.Range("D17").Value = 0.1626423157238523
$ws.Range("E17").Value = 0.1272325958003986
$ws.Range("F17").Value = 1.470455804104574
$ws.Range("G17").Value = 0.002411444506318998
$ws.Range("J17").Value = 0.1320212400530094
$ws.Range("M17").Value = 1.085782746330338
$ws.Range("N17").Value = 1.710604311494137
$ws.Range("O17").Value = 3.783184221272791
$ws.Range("B18").Value = 0.1734467614257369
$ws.Range("D18").Value = 0.162907490536373
$ws.Range("E18").Value = 0.1276168808607041
$ws.Range("F18").Value = 1.463406481868716
$ws.Range("G18").Value = 0.002412366332681941
$ws.Range("J18").Value = 0.1325598761090809
$ws.Range("M18").Value = 1.064259550421937
$ws.Range("N18").Value = 1.695796544514849
$ws.Range("O18").Value = 3.758506462251546
$ws.Range("B19").Value = 0.1725290952392413
$ws.Range("D19").Value = 0.1629981163096232
$ws.Range("E19").Value = 0.1277480800368922
$ws.Range("F19").Value = 1.461035670979513
$ws.Range("G19").Value = 0.002412680655949008
$ws.Range("J19").Value = 0.1327437529329041
$ws.Range("M19").Value = 1.056971072454573
$ws.Range("N19").Value = 1.690791414394226
$ws.Range("O19").Value = 3.750192594836278
$ws.Range("B20").Value = 0.1766629005796858
$ws.Range("D20").Value = 0.1625936377615353
$ws.Range("E20").Value = 0.1271619896164875
$ws.Range("F20").Value = 1.471768065069952
$ws.Range("G20").Value = 0.002411274945338409
$ws.Range("J20").Value = 0.1319222652984973
$ws.Range("M20").Value = 1.089765676215407
$ws.Range("N20").Value = 1.713348932133073
$ws.Range("O20").Value = 3.787771293982871
$ws.Range("B21").Value = 0.1906245300408216
$ws.Range("D21").Value = 0.1612931842326386
$ws.Range("E21").Value = 0.1252684701044471
$ws.Range("F21").Value = 1.508941138478377
$ws.Range("G21").Value = 0.002406705170106611
$ws.Range("J21").Value = 0.1292669960813582
$ws.Range("M21").Value = 1.199868436150567
$ws.Range("N21").Value = 1.789723719824792
$ws.Range("O21").Value = 3.916933297227274
$ws.Range("B22").Value = 0.1998021863980455
$ws.Range("D22").Value = 0.1604868410447349
$ws.Range("E22").Value = 0.1240874065027411
$ws.Range("F22").Value = 1.53406932126255
$ws.Range("G22").Value = 0.002403832238368761
$ws.Range("J22").Value = 0.1276099785669675
$ws.Range("M22").Value = 1.271757719244107
$ws.Range("N22").Value = 1.840064907550186
$ws.Range("O22").Value = 4.003522189238424
$ws.Range("B23").Value = 0.1948991324655935
$ws.Range("D23").Value = 0.1609132174445875
$ws.Range("E23").Value = 0.124712591853342
$ws.Range("F23").Value = 1.520581460567925
$ws.Range("G23").Value = 0.002405355207359294
$ws.Range("J23").Value = 0.1284871755598997
$ws.Range("M23").Value = 1.233395881224382
$ws.Range("N23").Value = 1.813158533964838
$ws.Range("O23").Value = 3.957108456873641
$ws.Range("B24").Value = 0.1764356155425304
$ws.Range("D24").Value = 0.1626156294339278
$ws.Range("E24").Value = 0.1271938904222596
$ws.Range("F24").Value = 1.471174512035574
$ws.Range("G24").Value = 0.002411351562677366
$ws.Range("J24").Value = 0.1319669837416946
$ws.Range("M24").Value = 1.08796504448415
$ws.Range("N24").Value = 1.712107956388905
$ws.Range("O24").Value = 3.78569675812895
$ws.Range("B25").Value = 0.1567546016999444
$ws.Range("D25").Value = 0.1646393940455688
$ws.Range("E25").Value = 0.1301123505123956
$ws.Range("F25").Value = 1.42143727785043
$ws.Range("G25").Value = 0.002418311044535545
$ws.Range("J25").Value = 0.136055425368669
$ws.Range("M25").Value = 0.9308462884999216
$ws.Range("N25").Value = 1.60500070849514
$ws.Range("O25").Value = 3.610066053435673
